$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$username = "['Assanhando Livros |', '@assanhandolivro']"
# Leading "'" forces these purely-numeric-looking values to be stored as
# text (shared strings) instead of being auto-coerced to numbers, matching
# the source data where Comments/Replys/Likes are text like "13"/"42"/"751".
$comments = "'13"
$replys = "'42"
$likes = "'751"
$views = "52,2 mil"
$text = "22 de ago Não tem nada melhor que livro com rico apaixonado e vou provar:   Ela foi traída pelo noivo e numa viagem para superar o chifre acaba conhecendo um cara rico num passeio (não fazendo ideia que ele seja) viram amigos e conforme os dias vão passando ele lhe oferece um acordo... + 13"

$ws.Range("A2").Value = $username
$ws.Range("B2").Value = $comments
$ws.Range("C2").Value = $replys
$ws.Range("D2").Value = $likes
$ws.Range("E2").Value = $views
$ws.Range("F2").Value = $text

$ws.Range("A3").Value = $username
$ws.Range("B3").Value = $comments
$ws.Range("C3").Value = $replys
$ws.Range("D3").Value = $likes
$ws.Range("E3").Value = $views
$ws.Range("F3").Value = $text

# Clear the implicit "quote prefix" style the apostrophe entry applied above
# so the cells keep the workbook's default (unstyled) formatting.
$ws.Range("B2:D3").Style = "Normal"
